$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.641.76'
$ws.Range("E2").Value = '  -1.93%  '

$ws.Range("D3").Value = '1.588.80'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("E4").Value = '  +0.19%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.92'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '

$ws.Range("E6").Value = '  -3.19%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  -2.69%  '

$ws.Range("E9").Value = '  -2.08%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.58'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.04%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0832'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.86%  '

$ws.Range("D12").Value = '1.809.82'
$ws.Range("E12").Value = '  -2.43%  '

$ws.Range("D13").Value = '1.578.45'
$ws.Range("E13").Value = '  -3.60%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.02'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("E15").Value = '  -4.58%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.68'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '26.613.11'
$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -2.50%  '

$ws.Range("E19").Value = '  +0.28%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '207.95'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.42%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.72'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.55%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.58%  '

$ws.Range("E23").Value = '  -4.06%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.52%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '146.62'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("E26").Value = '  +0.05%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.22'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.58%  '

$ws.Range("E28").Value = '  -3.94%  '

$ws.Range("E29").Value = '  -2.34%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0505'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.49%  '

$ws.Range("E31").Value = '  -2.00%  '

$ws.Range("E32").Value = '  -4.37%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.659'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +19.56%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.90'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.18%  '

$ws.Range("D35").Value = '1.303.90'
$ws.Range("E35").Value = '  -3.58%  '

$ws.Range("E36").Value = '  -1.17%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.48'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -5.67%  '

$ws.Range("E38").Value = '  -3.49%  '

$ws.Range("E39").Value = '  -3.49%  '

$ws.Range("E41").Value = '  -1.68%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.35'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '

$ws.Range("E43").Value = '  -3.41%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '62.50'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.55%  '

$ws.Range("D45").Value = '1.723.39'
$ws.Range("E45").Value = '  -2.27%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '89.43'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.43%  '

$ws.Range("E47").Value = '  -1.13%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.839'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.91%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0504'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.82%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0975'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.72%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.50'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.21%  '
